$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.140.00"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.634.84"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.17"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.96"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("E9").Value = "  -5.37%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.334"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "3.100.93"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "59.121.05"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.60"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.674.29"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "343.23"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.51"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.32"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.37"
$ws.Range("E23").Value = "  +4.07%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "2.759.17"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.11"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "0.0₃0788"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("E31").Value = "  -5.79%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.94"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.93"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.11"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.35"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.839"
$ws.Range("E38").Value = "  -6.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.836"
$ws.Range("E39").Value = "  -5.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0974"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.597"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.72"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "267.71"
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.07"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "2.030.42"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.64"
$ws.Range("E51").Value = "  -3.00%  "
